$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 'Faculty reviews'

$ws.Range("G3").Value = 'https://openalex.org/S4210207849'

$ws.Range("I3").Value = '2732-432X'

$ws.Range("D4").Value = 'Cell-based therapies potentially delay the trajectory toward end-stage kidney disease (ESKD) in late stage 4 diabetic chronic kidney disease (DKD). We describe the trial design, baseline patient characteristics, and early results of an IRB-approved phase II multicenter clinical trial, utilizing Renal Autologous Cell Therapy (REACT) in adults with pre-ESKD due to type 2 DKD. The trial objectives were safety and tolerability of REACT by assessment of the procedure, product administration, and renal-specific adverse events in addition to evaluate the impact on renal function following injection.Ten adults with an eGFR of 14-20 mL/min/1.73 m2 were enrolled in a single-arm open-label trial. Following a percutaneous kidney biopsy, an ex vivo expansion of selected renal cells that form the REACT was injected into the cortex of the biopsied kidney with CT image guidance. Each participant received two doses of the REACT product at 6-month intervals. A 6-month observation pre-trial was required to establish patients'' "own" baseline and rate of DKD progression.Five men and 5 women underwent 19 REACT injections (1 participant received only one injection). Their baseline characteristics were as follows: 3 Hispanic/Latino, 7 non-Hispanic, 7 white; mean (SD) age: 58.9 years (5.22), BMI 35.8 (8.2), eGFR (sCR) 15.5 (2.72), eGFR (sCR + Cys-C) 17.7 (3.67) mL/min/1.73 m2, sCr 3.6 mg/mL (0.73), Cys-C 2.6 mg/mL (0.52), and log random UACR 7.9 mg/g (1.01). The pre- and post-injection eGFR slope was -6.5 mL/min/1.73 m2 and -3.9 mL/min/1.73 m2. No cell-related adverse events occurred, and two procedure-related hematomas required observation without transfusion or angiographic interventions. Dialysis was delayed a mean of 16 months (range 6-28 months). At 15 months, 2 patients (20%) have eGFR slope stability and have not commenced renal replacement therapy.Trials that include patients with an eGFR of <20 mL/min/1.73 m2 are uncommon, and none to date involve autologous homologous cell-based treatments. REACT has the potential to stabilize or delay dialysis in high-risk late stage 4 DKD.'

$ws.Range("AF4").Value = 'c("https://openalex.org/W2503350049", "https://openalex.org/W2329386257", "https://openalex.org/W2397616145", "https://openalex.org/W1988477118", "https://openalex.org/W4324058133", "https://openalex.org/W2397320258", "https://openalex.org/W648618910", "https://openalex.org/W4296658948", "https://openalex.org/W3167355683", "https://openalex.org/W3211565028")'

$ws.Range("AF11").Value = 'c("https://openalex.org/W4206669628", "https://openalex.org/W3009669391", "https://openalex.org/W4205317059", "https://openalex.org/W3084498529", "https://openalex.org/W3036314732", "https://openalex.org/W4200329650", "https://openalex.org/W3127156785", "https://openalex.org/W3005417802", "https://openalex.org/W3171943759", "https://openalex.org/W4226296940")'
